$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 267, shifting existing rows 267-308 down to 268-309.
$ws.Rows(267).Insert()

# Populate the newly inserted row 267 with the new record.
$ws.Range("A267").Value = 11
$ws.Range("B267").Value = "Vega Monumental Concepción"
$ws.Range("C267").Value = "Bíobío"
$ws.Range("D267").Value = 45209
$ws.Range("E267").Value = 8
$ws.Range("F267").Value = "Fruta"
$ws.Range("G267").Value = 100108
$ws.Range("H267").Value = "Tropicales y subtropicales"
$ws.Range("I267").Value = 100108005
$ws.Range("J267").Value = "Piña"
$ws.Range("K267").Value = "Caramelo"
$ws.Range("L267").Value = "Primera"
$ws.Range("M267").Value = 100
$ws.Range("N267").Value = 19000
$ws.Range("O267").Value = 20000
$ws.Range("P267").Value = 19500
$ws.Range("Q267").Value = "$/caja 12 unidades"
$ws.Range("R267").Value = "Ecuador"
$ws.Range("S267").Value = 1625
$ws.Range("T267").Value = 12
